# Update benchmark: 2026-01-24 06:41:09 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - EFT - Şube (YKB F column gains a value)
$ws.Range("F3").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 4 - EFT - ATM
$ws.Range("F4").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 5 - EFT - Mobil
$ws.Range("F5").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 6 - DÜZENLİ EFT (AKBANK D column cleared)
$ws.Range("D6").Value = ""

# Row 8 - HESAPTAN HAVALE - Şube
$ws.Range("F8").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 9 - HESAPTAN HAVALE - ATM
$ws.Range("F9").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 10 - HESAPTAN HAVALE - Mobil
$ws.Range("F10").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 12 - GİDEN SWIFT (AKBANK D column cleared)
$ws.Range("D12").Value = ""

# Row 13 - GELEN SWIFT
$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 0,94 TL"
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 1.114 TL"
$ws.Range("F13").Value = "Hesaba: Asgari 300 TL | Azami 3.080 TL"

# Row 14 - GİDEN SWIFT - Mobil
$ws.Range("D14").Value = ""
$ws.Range("F14").Value = "1.952,38 TL - 9.523,81 TL"
